# NAT_holdings.xlsx update:
#   - Bump the "Model holdings provided as of" date in the confidential
#     disclaimer from 2021-03-18 to 2021-03-19.
#   - Refresh the Weight (D) / Percent Change (E) figures for rows 2-9.
#
# The worksheet ships protected, so we briefly unprotect it to make the
# writes, then restore protection afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Disclaimer text (cell A12): roll the as-of date forward a day ---
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."
# Writing the two-line string auto-expands the row; restore the original
# (non-custom) row height so row 12 stays exactly as it was.
$ws.Rows.Item(12).AutoFit()

# --- Weight / Percent Change figures ---
$ws.Range("D2").Value = 0.1496491938034782
$ws.Range("E2").Value = -0.0008857395925598865

$ws.Range("D3").Value = 0.1503119448832102
$ws.Range("E3").Value = 0

$ws.Range("D4").Value = 0.3006981982208146
$ws.Range("E4").Value = 0

$ws.Range("D5").Value = 0.149299743234165
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.07995489275984612
$ws.Range("E6").Value = 0

$ws.Range("D7").Value = 0.1202692376052996
$ws.Range("E7").Value = -0.0009852216748769127

$ws.Range("D8").Value = 0.04981678949318621
$ws.Range("E8").Value = 0

$ws.Range("E9").Value = -0.0002510420756560716

$ws.Protect()
